$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "CountryAndState" worksheet between "Login" and "Place" ---
$loginSheet = $wb.Worksheets.Item("Login")
$newSheet = $wb.Worksheets.Add($null, $loginSheet)
$newSheet.Name = "CountryAndState"

# --- 2. Populate the new sheet ---
$newSheet.Range("A1").Value = "Country"
$newSheet.Range("B1").Value = "State"
$newSheet.Range("A2").Value = "United States"
$newSheet.Range("B2").Value = "Texas"

# --- 3. Column width / autofit for the new sheet (column A only) ---
$newSheet.Columns.Item(1).ColumnWidth = 12.0221354166667

# --- 4. Update Login sheet's column widths (A & B) ---
$loginSheet.Columns.Item(1).ColumnWidth = 20.7369791666667
$loginSheet.Columns.Item(2).ColumnWidth = 12.7369791666667

# --- 5. Selection on Login moves to C2 ---
$loginSheet.Range("C2").Select() | Out-Null

# --- 6. New sheet becomes the active / selected tab, with its own selection at C9 ---
$newSheet.Range("C9").Select() | Out-Null
